$d = $word.ActiveDocument

# 1) Merge the three "D.2.2" / ". " / "Use Cases" runs into a single run.
$d.Content.Find.Execute("D.2.2. Use Cases", $true, $false, $false, $false, $false,
                         $true, 1, $false, "D.2.2. Use Cases", 2)

# 2) Move the "_GoBack" bookmark from the end of the document to the very
#    start (re-adding a bookmark with the same name relocates it).
$startRng = $d.Range(0, 1)
$d.Bookmarks.Add("_GoBack", $startRng)

# 3) Append the new "Group Participation" section (heading, intro
#    paragraph, spacer paragraph, participation table, spacer paragraph)
#    at the end of the document body, before the sectPr.
$endRng = $d.Content
$endRng.Collapse(0)
$newContentXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Heading1"/>
  </w:pPr>
  <w:r>
    <w:t>Group Participation</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t>Listed below is a table containing the group participation weights for each team member.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:tblPr>
    <w:tblStyle w:val="TableGrid"/>
    <w:tblW w:w="0" w:type="auto"/>
    <w:jc w:val="center"/>
    <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="2785"/>
    <w:gridCol w:w="2070"/>
  </w:tblGrid>
  <w:tr>
    <w:trPr>
      <w:jc w:val="center"/>
    </w:trPr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="2785" w:type="dxa"/>
        <w:shd w:val="clear" w:color="auto" w:fill="D0CECE" w:themeFill="background2" w:themeFillShade="E6"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:b/>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:b/>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
          <w:t>Team Member</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="2070" w:type="dxa"/>
        <w:shd w:val="clear" w:color="auto" w:fill="D0CECE" w:themeFill="background2" w:themeFillShade="E6"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:b/>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:b/>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
          <w:t>Participation</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:trPr>
      <w:jc w:val="center"/>
    </w:trPr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="2785" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
          <w:t>James Beasley</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="2070" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
          <w:t>25</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
          <w:t>%</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:trPr>
      <w:jc w:val="center"/>
    </w:trPr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="2785" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
          <w:t>Charles Beck</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="2070" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
          <w:t>25</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
          <w:t>%</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:trPr>
      <w:jc w:val="center"/>
    </w:trPr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="2785" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
          <w:t>Charles Duso</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="2070" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
          <w:t>25</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
          <w:t>%</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:trPr>
      <w:jc w:val="center"/>
    </w:trPr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="2785" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
          <w:t xml:space="preserve">Alexander </w:t>
        </w:r>
        <w:proofErr w:type="spellStart"/>
        <w:r>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
          <w:t>Grzesiak</w:t>
        </w:r>
        <w:proofErr w:type="spellEnd"/>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="2070" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
          <w:t>25</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
          <w:t>%</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:trPr>
      <w:jc w:val="center"/>
    </w:trPr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="2785" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
          <w:t>Erik Strauss</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="2070" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:keepNext/>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:sz w:val="32"/>
            <w:szCs w:val="32"/>
          </w:rPr>
          <w:t>0%</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
</w:tbl>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@
$endRng.InsertXML($newContentXml)

Write-Output "done"
